# Update existing rows 173-175 with revised figures from MV -datos- source
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 173 (Serie 184, 01-04-2021)
$ws.Cells.Item(173, 2).Value = 7637
$ws.Cells.Item(173, 3).Value = 7637
$ws.Cells.Item(173, 4).Value = 5890
$ws.Cells.Item(173, 9).Value = 648
$ws.Cells.Item(173, 10).Value = 286

# Row 174 (Serie 185, 01-05-2021)
$ws.Cells.Item(174, 2).Value = 6685
$ws.Cells.Item(174, 3).Value = 6685
$ws.Cells.Item(174, 4).Value = 4515
$ws.Cells.Item(174, 5).Value = 2170
$ws.Cells.Item(174, 8).Value = 2602
$ws.Cells.Item(174, 9).Value = 569
$ws.Cells.Item(174, 10).Value = 342

# Row 175 (Serie 186, 01-06-2021)
$ws.Cells.Item(175, 2).Value = 6991
$ws.Cells.Item(175, 3).Value = 6991
$ws.Cells.Item(175, 4).Value = 4908
$ws.Cells.Item(175, 7).Value = 2143
$ws.Cells.Item(175, 8).Value = 2740
$ws.Cells.Item(175, 9).Value = 718

# Add new row 176 (Serie 187, 01-07-2021)
# Force text format so the date-like label is stored as a shared string,
# matching column A's existing entries, then restore the default style.
$ws.Cells.Item(176, 1).NumberFormat = "@"
$ws.Cells.Item(176, 1).Value = "01-07-2021"
$ws.Cells.Item(176, 1).Style = "Normal"
$ws.Cells.Item(176, 2).Value = 6247
$ws.Cells.Item(176, 3).Value = 6247
$ws.Cells.Item(176, 4).Value = 4156
$ws.Cells.Item(176, 5).Value = 2092
$ws.Cells.Item(176, 6).Value = 638
$ws.Cells.Item(176, 7).Value = 1615
$ws.Cells.Item(176, 8).Value = 2703
$ws.Cells.Item(176, 9).Value = 632
$ws.Cells.Item(176, 10).Value = 299
$ws.Cells.Item(176, 11).Value = 113
$ws.Cells.Item(176, 12).Value = 248
$ws.Cells.Item(176, 13).Value = 0
$ws.Cells.Item(176, 14).Value = 0
$ws.Cells.Item(176, 15).Value = 0
